$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.239.98"
Set-TextValue $ws.Range("E2") "  -1.65%  "
Set-TextValue $ws.Range("D3") "1.671.76"
Set-TextValue $ws.Range("E3") "  +0.60%  "
Set-TextValue $ws.Range("E4") "  -0.18%  "
Set-TextValue $ws.Range("D5") "217.64"
Set-TextValue $ws.Range("E5") "  -1.07%  "
Set-TextValue $ws.Range("D6") "0.5128"
Set-TextValue $ws.Range("E6") "  +1.50%  "
Set-TextValue $ws.Range("E7") "  +0.12%  "
Set-TextValue $ws.Range("D8") "0.2658"
Set-TextValue $ws.Range("E8") "  +4.04%  "
Set-TextValue $ws.Range("D9") "0.06378"
Set-TextValue $ws.Range("E9") "  +3.69%  "
Set-TextValue $ws.Range("E10") "  -0.55%  "
Set-TextValue $ws.Range("D11") "0.07390"
Set-TextValue $ws.Range("E11") "  +0.66%  "
Set-TextValue $ws.Range("D12") "1.673.99"
Set-TextValue $ws.Range("E12") "  -0.67%  "
Set-TextValue $ws.Range("D13") "4.548"
Set-TextValue $ws.Range("E13") "  +2.22%  "
Set-TextValue $ws.Range("D14") "0.5835"
Set-TextValue $ws.Range("E14") "  +1.55%  "
Set-TextValue $ws.Range("D15") "1.900.53"
Set-TextValue $ws.Range("E15") "  +1.13%  "
Set-TextValue $ws.Range("D16") "0.000008657"
Set-TextValue $ws.Range("E16") "  +7.96%  "
Set-TextValue $ws.Range("D17") "64.58"
Set-TextValue $ws.Range("E17") "  +0.31%  "
Set-TextValue $ws.Range("D18") "26.297.95"
Set-TextValue $ws.Range("E18") "  -1.09%  "
Set-TextValue $ws.Range("D19") "4.958"
Set-TextValue $ws.Range("E19") "  +0.46%  "
Set-TextValue $ws.Range("E20") "  -0.08%  "
Set-TextValue $ws.Range("E21") "  +3.14%  "
Set-TextValue $ws.Range("D22") "189.25"
Set-TextValue $ws.Range("E22") "  +5.23%  "
Set-TextValue $ws.Range("D23") "6.212"
Set-TextValue $ws.Range("E23") "  +0.30%  "
Set-TextValue $ws.Range("E24") "  -0.13%  "
Set-TextValue $ws.Range("D25") "144.24"
Set-TextValue $ws.Range("E25") "  +1.06%  "
Set-TextValue $ws.Range("D26") "7.640"
Set-TextValue $ws.Range("E26") "  +1.49%  "
Set-TextValue $ws.Range("D27") "0.1182"
Set-TextValue $ws.Range("E27") "  +4.23%  "
Set-TextValue $ws.Range("E28") "  +4.73%  "
Set-TextValue $ws.Range("D29") "0.05969"
Set-TextValue $ws.Range("E29") "  +3.36%  "
Set-TextValue $ws.Range("D30") "1.283"
Set-TextValue $ws.Range("E30") "  -2.95%  "
Set-TextValue $ws.Range("E31") "  -1.01%  "
Set-TextValue $ws.Range("D32") "3.525"
Set-TextValue $ws.Range("E32") "  +3.18%  "
Set-TextValue $ws.Range("E33") "  +4.77%  "
Set-TextValue $ws.Range("E34") "  +3.09%  "
Set-TextValue $ws.Range("E35") "  +4.27%  "
Set-TextValue $ws.Range("D36") "0.6021"
Set-TextValue $ws.Range("E36") "  +0.84%  "
Set-TextValue $ws.Range("D37") "2.373"
Set-TextValue $ws.Range("E37") "  -1.97%  "
Set-TextValue $ws.Range("D38") "2.646"
Set-TextValue $ws.Range("E38") "  +0.55%  "
Set-TextValue $ws.Range("D39") "0.01620"
Set-TextValue $ws.Range("E39") "  +1.63%  "
Set-TextValue $ws.Range("D40") "6.082"
Set-TextValue $ws.Range("E40") "  +6.41%  "
Set-TextValue $ws.Range("D41") "1.081.98"
Set-TextValue $ws.Range("E41") "  +1.33%  "
Set-TextValue $ws.Range("D42") "0.8682"
Set-TextValue $ws.Range("E42") "  +1.29%  "
Set-TextValue $ws.Range("E43") "  +0.08%  "
Set-TextValue $ws.Range("D44") "100.21"
Set-TextValue $ws.Range("E44") "  +3.76%  "
Set-TextValue $ws.Range("E45") "  +1.52%  "
Set-TextValue $ws.Range("D46") "0.00000000115"
Set-TextValue $ws.Range("E46") "  +10.15%  "
Set-TextValue $ws.Range("D47") "56.14"
Set-TextValue $ws.Range("E47") "  +1.08%  "
Set-TextValue $ws.Range("D48") "1.010"
Set-TextValue $ws.Range("E48") "  +0.06%  "
Set-TextValue $ws.Range("D49") "8.082"
Set-TextValue $ws.Range("E49") "  +3.93%  "
Set-TextValue $ws.Range("E50") "  +0.14%  "
Set-TextValue $ws.Range("D51") "0.4295"
Set-TextValue $ws.Range("E51") "  -1.83%  "
